$d = $word.ActiveDocument

# The 3rd paragraph ("丰富单位的行为类（高）") and the 4th paragraph
# ("有效管理并复现关键帧（高）") need to become a single list item that only
# keeps the 4th paragraph's text, while preserving the "_GoBack" bookmark that
# sits inside the 3rd paragraph's runs (between the "（" and "高" runs).
#
# Strategy: merge the two paragraphs first (while the bookmark still anchors
# to real text, so it doesn't collapse to a zero-width range that a later
# Range.Delete() could sweep away), then trim the now-merged paragraph's text
# down to just the desired remainder, using the bookmark position to scope
# each deletion precisely (so we don't touch the identical "高）" substring
# that also occurs later in the document).

$p = $d.Paragraphs.Item(3)
$mark = $d.Range($p.Range.End - 1, $p.Range.End)
$mark.Delete()

$bm = $d.Bookmarks.Item("_GoBack")

$p = $d.Paragraphs.Item(3)
$before = $d.Range($p.Range.Start, $bm.Start)
$before.Delete()

$bm = $d.Bookmarks.Item("_GoBack")
$after = $d.Range($bm.End, $bm.End + 2)
$after.Delete()
